# Updated cryptos list values (prices / 1h volume deltas) and restored the
# Hedera / InjectiveProtocol, Kaspa / PEPE, Fetch.AI / Stellar row ordering
# that the upstream coinranking.com scrape produced on this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value for every cell the refreshed data touched
$updates = [ordered]@{
    "D2" = "72.002.91"
    "E2" = "  -0.67%  "
    "D3" = "3.998.05"
    "E3" = "  -1.09%  "
    "D4" = "1.00"
    "E4" = "  -0.19%  "
    "D5" = "542.96"
    "E5" = "  +4.18%  "
    "D6" = "150.46"
    "E6" = "  +1.53%  "
    "D7" = "0.702"
    "E7" = "  +11.89%  "
    "E8" = "  +0.14%  "
    "D9" = "0.747"
    "E9" = "  +0.76%  "
    "E10" = "  -3.17%  "
    "D11" = "53.01"
    "E11" = "  +11.54%  "
    "D12" = "0.0000325"
    "E12" = "  -3.01%  "
    "D13" = "10.67"
    "E13" = "  -2.27%  "
    "D14" = "4.637.40"
    "E14" = "  -1.00%  "
    "D15" = "3.987.81"
    "E15" = "  -1.91%  "
    "D16" = "14.16"
    "E16" = "  -0.61%  "
    "D17" = "20.56"
    "E17" = "  -3.32%  "
    "E18" = "  -0.20%  "
    "D19" = "1.19"
    "E19" = "  -2.05%  "
    "D20" = "71.849.34"
    "E20" = "  -0.92%  "
    "D21" = "432.91"
    "E21" = "  -1.94%  "
    "D22" = "96.89"
    "E22" = "  -4.51%  "
    "D23" = "3.55"
    "E23" = "  -0.62%  "
    "D24" = "4.28"
    "E24" = "  +5.93%  "
    "E25" = "  -2.46%  "
    "D26" = "11.58"
    "E26" = "  -2.99%  "
    "D27" = "10.75"
    "E27" = "  -4.61%  "
    "E28" = "  +1.06%  "
    "D29" = "36.78"
    "E29" = "  -2.62%  "
    "D30" = "3.61"
    "E30" = "  +16.83%  "
    "D31" = "7.48"
    "E31" = "  +7.73%  "
    "B32" = "Hedera"
    "C32" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D32" = "0.132"
    "E32" = "  +1.74%  "
    "E33" = "  -0.98%  "
    "B34" = "InjectiveProtocol"
    "C34" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D34" = "49.20"
    "E34" = "  +18.55%  "
    "D35" = "679.07"
    "E35" = "  -2.31%  "
    "D36" = "66.07"
    "E36" = "  -3.34%  "
    "D37" = "0.446"
    "E37" = "  +1.02%  "
    "B38" = "Kaspa"
    "C38" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D38" = "0.152"
    "E38" = "  -0.75%  "
    "B39" = "PEPE"
    "C39" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D39" = "0.0₃0831"
    "E39" = "  -7.14%  "
    "D40" = "3.41"
    "E40" = "  -7.46%  "
    "E41" = "  +5.54%  "
    "E42" = "  +0.08%  "
    "E43" = "  +0.22%  "
    "D44" = "0.0489"
    "E44" = "  -0.70%  "
    "B45" = "Fetch.AI"
    "C45" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D45" = "2.75"
    "E45" = "  -1.77%  "
    "B46" = "Stellar"
    "C46" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D46" = "0.150"
    "E46" = "  +2.03%  "
    "D47" = "9.85"
    "E47" = "  +7.86%  "
    "D48" = "3.36"
    "E48" = "  -3.99%  "
    "E49" = "  +1.13%  "
    "E50" = "  -4.11%  "
    "D51" = "144.80"
    "E51" = "  +1.30%  "
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    # Price/volume text such as "1.00", "542.96" or "0.150" must stay text -
    # otherwise Excel's auto-detection would coerce it to a number and drop
    # the significant trailing/leading zeros. Force text, write, then drop
    # back to the default "Normal" style so no stray number format sticks.
    $looksNumeric = $value -match '^[+-]?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

